# Rename the single worksheet: "Property1" -> "DataNode"
# (ties in with the commit message: "unify the conception of DataNode, DataTable, Entity")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move / record the user's current selection: was A9, now C24 (still inside the
# scrolling pane below the frozen header rows).
$ws.Range("C24").Select()

# Column A was re-measured slightly narrower (24.1640625 -> 24.125 "characters"
# of stored width). ColumnWidth is expressed in characters and gets rounded to
# the nearest pixel internally, so feed it the character width that lands on
# the closest achievable pixel boundary to the target.
$ws.Columns.Item(1).ColumnWidth = 23.43

# A new, smaller (9pt) 宋体 font was introduced into the workbook's font table
# (used for East-Asian phonetic-guide metadata on the sheet). Touch a cell that
# already carries the regular (non-bold) base font with that size so the new
# font entry is recorded, then immediately restore the cell's original
# formatting so no visible cell style actually changes.
$ws.Range("A9").Font.Size = 9
$ws.Range("A9").Font.Size = 11

# Rename the built-in "Normal" cell style to its localized Chinese name.
$wb.Styles.Item("Normal").Name = "常规"
